$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 68.158272
$ws.Range("H2").Value = 204.474816
$ws.Range("I2").Value = 0.164824640128582
$ws.Range("J2").Value = 0.1648246401285819
$ws.Range("M2").Value = 10.13412066666667
$ws.Range("N2").Value = 30.402362
$ws.Range("O2").Value = 0.3332793623493037
$ws.Range("P2").Value = 0.3332793623493037
$ws.Range("Q2").Value = 690.724152879488
$ws.Range("R2").Value = 6216.517375915392
$ws.Range("S2").Value = 0.05493265096150726
$ws.Range("T2").Value = 0.05493265096150724

# Row 3
$ws.Range("G3").Value = 68.158272
$ws.Range("H3").Value = 204.474816
$ws.Range("I3").Value = 0.164824640128582
$ws.Range("J3").Value = 0.1648246401285819
$ws.Range("O3").Value = 0.2828502524658126
$ws.Range("P3").Value = 0.2828502524658126
$ws.Range("Q3").Value = 586.2094179759999
$ws.Range("R3").Value = 5275.884761784
$ws.Range("S3").Value = 0.04662069107295613
$ws.Range("T3").Value = 0.04662069107295611

# Row 4
$ws.Range("G4").Value = 68.158272
$ws.Range("H4").Value = 204.474816
$ws.Range("I4").Value = 0.164824640128582
$ws.Range("J4").Value = 0.1648246401285819
$ws.Range("M4").Value = 11.67245633333333
$ws.Range("N4").Value = 35.017369
$ws.Range("O4").Value = 0.3838703851848838
$ws.Range("P4").Value = 0.3838703851848837
$ws.Range("Q4").Value = 795.5744536754559
$ws.Range("R4").Value = 7160.170083079103
$ws.Range("S4").Value = 0.06327129809411861
$ws.Range("T4").Value = 0.06327129809411859

# Row 5
$ws.Range("I5").Value = 0.3471155005059974
$ws.Range("J5").Value = 0.3471155005059974
$ws.Range("M5").Value = 10.13412066666667
$ws.Range("N5").Value = 30.402362
$ws.Range("O5").Value = 0.3332793623493037
$ws.Range("P5").Value = 0.3332793623493037
$ws.Range("Q5").Value = 1454.643309709663
$ws.Range("R5").Value = 13091.78978738696
$ws.Range("S5").Value = 0.1156864326701982
$ws.Range("T5").Value = 0.1156864326701982

# Row 6
$ws.Range("I6").Value = 0.3471155005059974
$ws.Range("J6").Value = 0.3471155005059974
$ws.Range("O6").Value = 0.2828502524658126
$ws.Range("P6").Value = 0.2828502524658126
$ws.Range("S6").Value = 0.09818170695291829
$ws.Range("T6").Value = 0.09818170695291825

# Row 7
$ws.Range("I7").Value = 0.3471155005059974
$ws.Range("J7").Value = 0.3471155005059974
$ws.Range("M7").Value = 11.67245633333333
$ws.Range("N7").Value = 35.017369
$ws.Range("O7").Value = 0.3838703851848838
$ws.Range("P7").Value = 0.3838703851848837
$ws.Range("Q7").Value = 1675.454740637735
$ws.Range("R7").Value = 15079.09266573962
$ws.Range("S7").Value = 0.1332473608828809
$ws.Range("T7").Value = 0.1332473608828809

# Row 8
$ws.Range("G8").Value = 201.822474
$ws.Range("H8").Value = 605.4674219999999
$ws.Range("I8").Value = 0.4880598593654206
$ws.Range("J8").Value = 0.4880598593654206
$ws.Range("M8").Value = 10.13412066666667
$ws.Range("N8").Value = 30.402362
$ws.Range("O8").Value = 0.3332793623493037
$ws.Range("P8").Value = 0.3332793623493037
$ws.Range("Q8").Value = 2045.293304761195
$ws.Range("R8").Value = 18407.63974285076
$ws.Range("S8").Value = 0.1626602787175982
$ws.Range("T8").Value = 0.1626602787175982

# Row 9
$ws.Range("G9").Value = 201.822474
$ws.Range("H9").Value = 605.4674219999999
$ws.Range("I9").Value = 0.4880598593654206
$ws.Range("J9").Value = 0.4880598593654206
$ws.Range("O9").Value = 0.2828502524658126
$ws.Range("P9").Value = 0.2828502524658126
$ws.Range("Q9").Value = 1735.81623398575
$ws.Range("R9").Value = 15622.34610587175
$ws.Range("S9").Value = 0.1380478544399382
$ws.Range("T9").Value = 0.1380478544399382

# Row 10
$ws.Range("G10").Value = 201.822474
$ws.Range("H10").Value = 605.4674219999999
$ws.Range("I10").Value = 0.4880598593654206
$ws.Range("J10").Value = 0.4880598593654206
$ws.Range("M10").Value = 11.67245633333333
$ws.Range("N10").Value = 35.017369
$ws.Range("O10").Value = 0.3838703851848838
$ws.Range("P10").Value = 0.3838703851848837
$ws.Range("Q10").Value = 2355.764014850301
$ws.Range("R10").Value = 21201.87613365271
$ws.Range("S10").Value = 0.1873517262078842
$ws.Range("T10").Value = 0.1873517262078842
